# Update expense amounts in column E for the "set of counteragents added" change.
# Each row's total amount (column E) is increased to reflect the updated set.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E6").Value = 4380
$ws.Range("E8").Value = 7831
$ws.Range("E11").Value = 1256
$ws.Range("E12").Value = 1308
$ws.Range("E13").Value = 1308
$ws.Range("E14").Value = 1308
$ws.Range("E15").Value = 880
$ws.Range("E16").Value = 849
$ws.Range("E17").Value = 849
$ws.Range("E18").Value = 849
$ws.Range("E19").Value = 849
$ws.Range("E20").Value = 849
$ws.Range("E21").Value = 849
$ws.Range("E22").Value = 849
$ws.Range("E24").Value = 2761
$ws.Range("E25").Value = 2903
$ws.Range("E26").Value = 2605
$ws.Range("E27").Value = 2523
$ws.Range("E28").Value = 2523
$ws.Range("E29").Value = 2523
$ws.Range("E30").Value = 2523
$ws.Range("E31").Value = 2523
$ws.Range("E32").Value = 2523
$ws.Range("E33").Value = 2523
$ws.Range("E34").Value = 2523
$ws.Range("E35").Value = 2523
$ws.Range("E40").Value = 56000
$ws.Range("E41").Value = 597545
$ws.Range("E46").Value = 186600
$ws.Range("E47").Value = 60550
$ws.Range("E71").Value = 203962
$ws.Range("E79").Value = 11500
$ws.Range("E80").Value = 6930
$ws.Range("E99").Value = 11900
$ws.Range("E100").Value = 26833
$ws.Range("E101").Value = 302461
$ws.Range("E103").Value = 15500
$ws.Range("E104").Value = 26557
$ws.Range("E108").Value = 21671
$ws.Range("E112").Value = 45000
$ws.Range("E114").Value = 4480
$ws.Range("E116").Value = 12600
$ws.Range("E118").Value = 135000
$ws.Range("E119").Value = 8500
